$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (ECs -> ECs self-signaling row), changing Target cluster and numeric values
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.02506566666666667
$ws.Range("H2").Value = 0.075197
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.52656
$ws.Range("N2").Value = 1.57968
$ws.Range("O2").Value = 0.9686255056421601
$ws.Range("P2").Value = 0.9686255056421602
$ws.Range("Q2").Value = 0.01319857744
$ws.Range("R2").Value = 0.11878719696
$ws.Range("S2").Value = 0.9686255056421601
$ws.Range("T2").Value = 0.9686255056421602

# Add new row 3 (ECs -> FAPs), which holds what used to be row 2's values
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dkk4"
$ws.Range("C3").Value = "Kremen2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02506566666666667
$ws.Range("H3").Value = 0.075197
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01705566666666667
$ws.Range("N3").Value = 0.051167
$ws.Range("O3").Value = 0.03137449435783982
$ws.Range("P3").Value = 0.03137449435783982
$ws.Range("Q3").Value = 0.0004275116554444444
$ws.Range("R3").Value = 0.003847604899
$ws.Range("S3").Value = 0.03137449435783982
$ws.Range("T3").Value = 0.03137449435783982
